# Update database and change read_price algorithm:
# Rows 11-27, columns D:H on the "Overview" sheet hold the yearly figures.
# Under the new read_price algorithm every figure resets to 0, except the
# two "not applicable" rows (15 and 23), which use the "-" placeholder
# text already used elsewhere in the sheet (e.g. D15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$dashRows = @(15, 23)
$zeroRows = @(11, 12, 13, 14, 16, 17, 18, 19, 20, 21, 22, 24, 25, 26, 27)
$cols = @("D", "E", "F", "G", "H")

foreach ($r in $zeroRows) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = 0
    }
}

foreach ($r in $dashRows) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = "-"
    }
}
